{"js": "// Update the date and all two-digit-by-two-digit multiplication\n// problems to the values from the new day's worksheet.\nconst replacements = [\n  [\"2025-10-23 Thursday\", \"2025-10-24 Friday\"],\n  [\"36\u00d723=\", \"59\u00d718=\"],\n  [\"25\u00d793=\", \"81\u00d747=\"],\n  [\"52\u00d775=\", \"37\u00d735=\"],\n  [\"83\u00d741=\", \"38\u00d782=\"],\n  [\"31\u00d736=\", \"38\u00d731=\"],\n  [\"80\u00d720=\", \"83\u00d744=\"],\n  [\"46\u00d781=\", \"87\u00d714=\"],\n  [\"18\u00d728=\", \"12\u00d722=\"],\n  [\"18\u00d767=\", \"84\u00d716=\"],\n  [\"72\u00d773=\", \"98\u00d796=\"],\n  [\"71\u00d797=\", \"71\u00d755=\"],\n  [\"92\u00d737=\", \"30\u00d729=\"],\n  [\"29\u00d719=\", \"65\u00d749=\"],\n  [\"60\u00d776=\", \"76\u00d727=\"],\n  [\"48\u00d753=\", \"55\u00d765=\"],\n  [\"75\u00d750=\", \"72\u00d785=\"],\n  [\"73\u00d737=\", \"15\u00d791=\"],\n  [\"23\u00d730=\", \"80\u00d798=\"],\n  [\"42\u00d764=\", \"48\u00d723=\"],\n  [\"65\u00d765=\", \"29\u00d728=\"],\n  [\"99\u00d789=\", \"68\u00d799=\"],\n  [\"84\u00d754=\", \"30\u00d782=\"],\n  [\"88\u00d737=\", \"26\u00d754=\"],\n  [\"65\u00d734=\", \"47\u00d745=\"],\n  [\"25\u00d714=\", \"81\u00d782=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date and all two-digit-by-two-digit multiplication\n# problems to the values from the new day's worksheet.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-10-23 Thursday\", \"2025-10-24 Friday\"),\n    @(\"36\u00d723=\", \"59\u00d718=\"),\n    @(\"25\u00d793=\", \"81\u00d747=\"),\n    @(\"52\u00d775=\", \"37\u00d735=\"),\n    @(\"83\u00d741=\", \"38\u00d782=\"),\n    @(\"31\u00d736=\", \"38\u00d731=\"),\n    @(\"80\u00d720=\", \"83\u00d744=\"),\n    @(\"46\u00d781=\", \"87\u00d714=\"),\n    @(\"18\u00d728=\", \"12\u00d722=\"),\n    @(\"18\u00d767=\", \"84\u00d716=\"),\n    @(\"72\u00d773=\", \"98\u00d796=\"),\n    @(\"71\u00d797=\", \"71\u00d755=\"),\n    @(\"92\u00d737=\", \"30\u00d729=\"),\n    @(\"29\u00d719=\", \"65\u00d749=\"),\n    @(\"60\u00d776=\", \"76\u00d727=\"),\n    @(\"48\u00d753=\", \"55\u00d765=\"),\n    @(\"75\u00d750=\", \"72\u00d785=\"),\n    @(\"73\u00d737=\", \"15\u00d791=\"),\n    @(\"23\u00d730=\", \"80\u00d798=\"),\n    @(\"42\u00d764=\", \"48\u00d723=\"),\n    @(\"65\u00d765=\", \"29\u00d728=\"),\n    @(\"99\u00d789=\", \"68\u00d799=\"),\n    @(\"84\u00d754=\", \"30\u00d782=\"),\n    @(\"88\u00d737=\", \"26\u00d754=\"),\n    @(\"65\u00d734=\", \"47\u00d745=\"),\n    @(\"25\u00d714=\", \"81\u00d782=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
